$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 111, pushing the existing row 111..151 data down to 112..152
$ws.Rows.Item(111).Insert()

# Populate the newly inserted row 111 with the new weekly record
$ws.Range("A111").Value = 4
$ws.Range("B111").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C111").Value = "Los Lagos"
$ws.Range("D111").Value = 45146
$ws.Range("E111").Value = 10
$ws.Range("F111").Value = 100112026
$ws.Range("G111").Value = "Haba"
$ws.Range("H111").Value = "Sin especificar"
$ws.Range("I111").Value = "Primera"
$ws.Range("J111").Value = 80
$ws.Range("K111").Value = 18000
$ws.Range("L111").Value = 18000
$ws.Range("M111").Value = 18000
$ws.Range("N111").Value = "$/saco 25 kilos"
$ws.Range("O111").Value = "Provincia de Limarí"
$ws.Range("P111").Value = 720
$ws.Range("Q111").Value = 25
$ws.Range("R111").Value = "Hortaliza"
